$d = $word.ActiveDocument

function Get-ParagraphAt($rng) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Start -le $rng.Start -and $p.Range.End -ge $rng.End) {
            return $p
        }
    }
    return $null
}

# --- 1. Executive Summary paragraph: merge runs (no text change) ---
$d.Content.Find.Execute(
    "The usability test of the website will be focused on all core functionality aspects of the system; UI elements are mostly done, but have not been fully implemented across the website. This includes logging in, logging out, navigation between pages, accessing the list of events, creating an event, and booking an event. As the calendar and payment functionalities are not tied to any databases, successful navigation to and from those pages will be the only tests performed for those pages. Any issues reported will be documented to ensure they are no longer a factor upon the final release of the Event planner.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The usability test of the website will be focused on all core functionality aspects of the system; UI elements are mostly done, but have not been fully implemented across the website. This includes logging in, logging out, navigation between pages, accessing the list of events, creating an event, and booking an event. As the calendar and payment functionalities are not tied to any databases, successful navigation to and from those pages will be the only tests performed for those pages. Any issues reported will be documented to ensure they are no longer a factor upon the final release of the Event planner.",
    2) | Out-Null

# --- 2. Methodology paragraph: merge runs (no text change) ---
$d.Content.Find.Execute(
    "We expect five students to handle the usability testing of this system, and will perform this test during class time.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "We expect five students to handle the usability testing of this system, and will perform this test during class time.",
    2) | Out-Null

# --- 3. Participants paragraph: merge runs + grammar fix "in to" -> "into" ---
$d.Content.Find.Execute(
    "We will have a target of five students to come to our demonstration setup during the testing session in class; given the average skill set of the class members, no additional characteristics will be necessary to have. The recruitment of testers will be handled through the class testing day, as all projects will be readily available for voluntary user testing. These participants will already have well beyond the minimum necessary background knowledge to properly use the program, and may even have css, html or javascript training that could be factored in to the feedback of the test. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "We will have a target of five students to come to our demonstration setup during the testing session in class; given the average skill set of the class members, no additional characteristics will be necessary to have. The recruitment of testers will be handled through the class testing day, as all projects will be readily available for voluntary user testing. These participants will already have well beyond the minimum necessary background knowledge to properly use the program, and may even have css, html or javascript training that could be factored into the feedback of the test. ",
    2) | Out-Null

# --- 4. Empty run right after the Participants paragraph: give it rPr (red, Verdana 10pt) ---
$rng4 = $d.Content
$rng4.Find.Execute("factored into the feedback of the test. ") | Out-Null
$para4 = Get-ParagraphAt $rng4
$afterPara4 = $para4.Next()
$afterPara4.Range.Font.Name = "Verdana"
$afterPara4.Range.Font.NameAscii = "Verdana"
$afterPara4.Range.Font.NameFarEast = "Verdana"
$afterPara4.Range.Font.NameBi = "Verdana"
$afterPara4.Range.Font.Color = 255
$afterPara4.Range.Font.Size = 10
$afterPara4.Range.Font.SizeBi = 10

# --- 5. "Due to the nature of the test environment..." paragraph: merge runs (no text change) ---
$d.Content.Find.Execute(
    "Due to the nature of the test environment, the participants will select our project for testing, rather than the other way around. Despite the less controlled nature of the participant selection, they will be more than capable of performing the tasks that are provided to them given the background required to enter CSCE 3444. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Due to the nature of the test environment, the participants will select our project for testing, rather than the other way around. Despite the less controlled nature of the participant selection, they will be more than capable of performing the tasks that are provided to them given the background required to enter CSCE 3444. ",
    2) | Out-Null

# --- 6. "As the system is a basic website..." run group: merge runs (no text change) ---
$d.Content.Find.Execute(
    "As the system is a basic website, a full training session would be unnecessary; however, our trainer will be guiding testers through the initial process of using the website to streamline the testing experience.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "As the system is a basic website, a full training session would be unnecessary; however, our trainer will be guiding testers through the initial process of using the website to streamline the testing experience.",
    2) | Out-Null

# --- 7. "The participants will receive an overview..." run group: merge runs (no text change) ---
$d.Content.Find.Execute(
    "The participants will receive an overview of the usability test procedure, equipment and software. Currently, the calendar and payment functions are properly displayed and can be navigated to, but do not connect to anything as of the time of testing.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The participants will receive an overview of the usability test procedure, equipment and software. Currently, the calendar and payment functions are properly displayed and can be navigated to, but do not connect to anything as of the time of testing.",
    2) | Out-Null

# --- 8. "Participants will take part in the usability test..." paragraph: merge runs (no text change) ---
$d.Content.Find.Execute(
    "Participants will take part in the usability test in room 005 of the Business Leadership Building. A Laptop with the Web site/Web application and supporting software will be used in a typical office environment.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Participants will take part in the usability test in room 005 of the Business Leadership Building. A Laptop with the Web site/Web application and supporting software will be used in a typical office environment.",
    2) | Out-Null

# --- 9. Second empty run (near the end, right after the "base process" paragraph): give it rPr ---
$rng9 = $d.Content
$rng9.Find.Execute("base process") | Out-Null
$para9 = Get-ParagraphAt $rng9
$afterPara9 = $para9.Next()
$afterPara9.Range.Font.Name = "Verdana"
$afterPara9.Range.Font.NameAscii = "Verdana"
$afterPara9.Range.Font.NameFarEast = "Verdana"
$afterPara9.Range.Font.NameBi = "Verdana"
$afterPara9.Range.Font.Color = 0
$afterPara9.Range.Font.Size = 10
$afterPara9.Range.Font.SizeBi = 10

# --- 10. Grammar fix: "and action" -> "an action" ---
$d.Content.Find.Execute(
    "initiates (or attempts to initiate) and action that will result",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "initiates (or attempts to initiate) an action that will result",
    2) | Out-Null
